# Weekly update: insert the newest week's two quality-grade rows
# (Primera / Segunda) at the top of the data block for
# "Agrícola del Norte S.A. de Arica - Pepino ensalada", pushing the
# existing historical rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the first data row of the
# block (row 207), shifting rows 207:310 down to 209:312.
$ws.Range("A207:A208").EntireRow.Insert()

# --- New row 207: "Primera" quality for the new week ---
$ws.Cells.Item(207, 1).Value = 1
$ws.Cells.Item(207, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(207, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(207, 4).Value = 44704
$ws.Cells.Item(207, 5).Value = 15
$ws.Cells.Item(207, 6).Value = 100112043
$ws.Cells.Item(207, 7).Value = "Pepino ensalada"
$ws.Cells.Item(207, 8).Value = "Sin especificar"
$ws.Cells.Item(207, 9).Value = "Primera"
$ws.Cells.Item(207, 10).Value = 130
$ws.Cells.Item(207, 11).Value = 16000
$ws.Cells.Item(207, 12).Value = 17000
$ws.Cells.Item(207, 13).Value = 16500
$ws.Cells.Item(207, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(207, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(207, 16).Value = 236
$ws.Cells.Item(207, 17).Value = 70
$ws.Cells.Item(207, 18).Value = "Hortaliza"

# --- New row 208: "Segunda" quality for the new week ---
$ws.Cells.Item(208, 1).Value = 1
$ws.Cells.Item(208, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(208, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(208, 4).Value = 44704
$ws.Cells.Item(208, 5).Value = 15
$ws.Cells.Item(208, 6).Value = 100112043
$ws.Cells.Item(208, 7).Value = "Pepino ensalada"
$ws.Cells.Item(208, 8).Value = "Sin especificar"
$ws.Cells.Item(208, 9).Value = "Segunda"
$ws.Cells.Item(208, 10).Value = 160
$ws.Cells.Item(208, 11).Value = 13000
$ws.Cells.Item(208, 12).Value = 14000
$ws.Cells.Item(208, 13).Value = 13500
$ws.Cells.Item(208, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(208, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(208, 16).Value = 135
$ws.Cells.Item(208, 17).Value = 100
$ws.Cells.Item(208, 18).Value = "Hortaliza"
